$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'Dcn'
$ws.Cells.Item(2,3).Value = 'Tlr2'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 12.265976
$ws.Cells.Item(2,8).Value = 36.797928
$ws.Cells.Item(2,9).Value = 0.004000867643088759
$ws.Cells.Item(2,10).Value = 0.004000867643088758
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.597802666666666
$ws.Cells.Item(2,14).Value = 4.793407999999999
$ws.Cells.Item(2,15).Value = 0.02304920886321625
$ws.Cells.Item(2,16).Value = 0.02304920886321625
$ws.Cells.Item(2,17).Value = 19.59860916206933
$ws.Cells.Item(2,18).Value = 176.387482458624
$ws.Cells.Item(2,19).Value = [double]"9.221683393963653E-05"
$ws.Cells.Item(2,20).Value = [double]"9.221683393963652E-05"

$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'Dcn'
$ws.Cells.Item(3,3).Value = 'Tlr2'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 12.265976
$ws.Cells.Item(3,8).Value = 36.797928
$ws.Cells.Item(3,9).Value = 0.004000867643088759
$ws.Cells.Item(3,10).Value = 0.004000867643088758
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 1.27306
$ws.Cells.Item(3,14).Value = 3.81918
$ws.Cells.Item(3,15).Value = 0.0183646118807784
$ws.Cells.Item(3,16).Value = 0.0183646118807784
$ws.Cells.Item(3,17).Value = 15.61532340656
$ws.Cells.Item(3,18).Value = 140.53791065904
$ws.Cells.Item(3,19).Value = [double]"7.347438145168971E-05"
$ws.Cells.Item(3,20).Value = [double]"7.34743814516897E-05"

$ws.Cells.Item(4,1).Value = 'ECs'
$ws.Cells.Item(4,2).Value = 'Dcn'
$ws.Cells.Item(4,3).Value = 'Tlr2'
$ws.Cells.Item(4,4).Value = 'Resolving-Mac'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 12.265976
$ws.Cells.Item(4,8).Value = 36.797928
$ws.Cells.Item(4,9).Value = 0.004000867643088759
$ws.Cells.Item(4,10).Value = 0.004000867643088758
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 66.45050433333334
$ws.Cells.Item(4,14).Value = 199.351513
$ws.Cells.Item(4,15).Value = 0.9585861792560053
$ws.Cells.Item(4,16).Value = 0.9585861792560054
$ws.Cells.Item(4,17).Value = 815.0802913405628
$ws.Cells.Item(4,18).Value = 7335.722622065065
$ws.Cells.Item(4,19).Value = 0.003835176427697433
$ws.Cells.Item(4,20).Value = 0.003835176427697433

$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'Dcn'
$ws.Cells.Item(5,3).Value = 'Tlr2'
$ws.Cells.Item(5,4).Value = 'ECs'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3042.696044666667
$ws.Cells.Item(5,8).Value = 9128.088134
$ws.Cells.Item(5,9).Value = 0.9924545876219728
$ws.Cells.Item(5,10).Value = 0.9924545876219727
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 1.597802666666666
$ws.Cells.Item(5,14).Value = 4.793407999999999
$ws.Cells.Item(5,15).Value = 0.02304920886321625
$ws.Cells.Item(5,16).Value = 0.02304920886321625
$ws.Cells.Item(5,17).Value = 4861.627854024518
$ws.Cells.Item(5,18).Value = 43754.65068622067
$ws.Cells.Item(5,19).Value = 0.022875293077356
$ws.Cells.Item(5,20).Value = 0.022875293077356

$ws.Cells.Item(6,1).Value = 'FAPs'
$ws.Cells.Item(6,2).Value = 'Dcn'
$ws.Cells.Item(6,3).Value = 'Tlr2'
$ws.Cells.Item(6,4).Value = 'FAPs'
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3042.696044666667
$ws.Cells.Item(6,8).Value = 9128.088134
$ws.Cells.Item(6,9).Value = 0.9924545876219728
$ws.Cells.Item(6,10).Value = 0.9924545876219727
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 1.27306
$ws.Cells.Item(6,14).Value = 3.81918
$ws.Cells.Item(6,15).Value = 0.0183646118807784
$ws.Cells.Item(6,16).Value = 0.0183646118807784
$ws.Cells.Item(6,17).Value = 3873.534626623346
$ws.Cells.Item(6,18).Value = 34861.81163961012
$ws.Cells.Item(6,19).Value = 0.01822604331097551
$ws.Cells.Item(6,20).Value = 0.01822604331097551

$ws.Cells.Item(7,1).Value = 'FAPs'
$ws.Cells.Item(7,2).Value = 'Dcn'
$ws.Cells.Item(7,3).Value = 'Tlr2'
$ws.Cells.Item(7,4).Value = 'Resolving-Mac'
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3042.696044666667
$ws.Cells.Item(7,8).Value = 9128.088134
$ws.Cells.Item(7,9).Value = 0.9924545876219728
$ws.Cells.Item(7,10).Value = 0.9924545876219727
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 66.45050433333334
$ws.Cells.Item(7,14).Value = 199.351513
$ws.Cells.Item(7,15).Value = 0.9585861792560053
$ws.Cells.Item(7,16).Value = 0.9585861792560054
$ws.Cells.Item(7,17).Value = 202188.6867011386
$ws.Cells.Item(7,18).Value = 1819698.180310247
$ws.Cells.Item(7,19).Value = 0.9513532512336412
$ws.Cells.Item(7,20).Value = 0.9513532512336412

$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Dcn'
$ws.Cells.Item(8,3).Value = 'Tlr2'
$ws.Cells.Item(8,4).Value = 'ECs'
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 8.377189333333332
$ws.Cells.Item(8,8).Value = 25.131568
$ws.Cells.Item(8,9).Value = 0.002732438555542716
$ws.Cells.Item(8,10).Value = 0.002732438555542716
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 1.597802666666666
$ws.Cells.Item(8,14).Value = 4.793407999999999
$ws.Cells.Item(8,15).Value = 0.02304920886321625
$ws.Cells.Item(8,16).Value = 0.02304920886321625
$ws.Cells.Item(8,17).Value = 13.38509545597155
$ws.Cells.Item(8,18).Value = 120.465859103744
$ws.Cells.Item(8,19).Value = [double]"6.298054697260897E-05"
$ws.Cells.Item(8,20).Value = [double]"6.298054697260897E-05"

$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Dcn'
$ws.Cells.Item(9,3).Value = 'Tlr2'
$ws.Cells.Item(9,4).Value = 'FAPs'
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 8.377189333333332
$ws.Cells.Item(9,8).Value = 25.131568
$ws.Cells.Item(9,9).Value = 0.002732438555542716
$ws.Cells.Item(9,10).Value = 0.002732438555542716
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 1.27306
$ws.Cells.Item(9,14).Value = 3.81918
$ws.Cells.Item(9,15).Value = 0.0183646118807784
$ws.Cells.Item(9,16).Value = 0.0183646118807784
$ws.Cells.Item(9,17).Value = 10.66466465269333
$ws.Cells.Item(9,18).Value = 95.98198187423999
$ws.Cells.Item(9,19).Value = [double]"5.018017356061673E-05"
$ws.Cells.Item(9,20).Value = [double]"5.018017356061674E-05"

$ws.Cells.Item(10,1).Value = 'MuSCs'
$ws.Cells.Item(10,2).Value = 'Dcn'
$ws.Cells.Item(10,3).Value = 'Tlr2'
$ws.Cells.Item(10,4).Value = 'Resolving-Mac'
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 8.377189333333332
$ws.Cells.Item(10,8).Value = 25.131568
$ws.Cells.Item(10,9).Value = 0.002732438555542716
$ws.Cells.Item(10,10).Value = 0.002732438555542716
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 66.45050433333334
$ws.Cells.Item(10,14).Value = 199.351513
$ws.Cells.Item(10,15).Value = 0.9585861792560053
$ws.Cells.Item(10,16).Value = 0.9585861792560054
$ws.Cells.Item(10,17).Value = 556.6684560958204
$ws.Cells.Item(10,18).Value = 5010.016104862384
$ws.Cells.Item(10,19).Value = 0.00261927783500949
$ws.Cells.Item(10,20).Value = 0.00261927783500949

$ws.Cells.Item(11,1).Value = 'Resolving-Mac'
$ws.Cells.Item(11,2).Value = 'Dcn'
$ws.Cells.Item(11,3).Value = 'Tlr2'
$ws.Cells.Item(11,4).Value = 'ECs'
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.489778666666667
$ws.Cells.Item(11,8).Value = 7.469336
$ws.Cells.Item(11,9).Value = 0.0008121061793956991
$ws.Cells.Item(11,10).Value = 0.0008121061793956991
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 1.597802666666666
$ws.Cells.Item(11,14).Value = 4.793407999999999
$ws.Cells.Item(11,15).Value = 0.02304920886321625
$ws.Cells.Item(11,16).Value = 0.02304920886321625
$ws.Cells.Item(11,17).Value = 3.978174993009777
$ws.Cells.Item(11,18).Value = 35.80357493708799
$ws.Cells.Item(11,19).Value = [double]"1.871840494800003E-05"
$ws.Cells.Item(11,20).Value = [double]"1.871840494800003E-05"

$ws.Cells.Item(12,1).Value = 'Resolving-Mac'
$ws.Cells.Item(12,2).Value = 'Dcn'
$ws.Cells.Item(12,3).Value = 'Tlr2'
$ws.Cells.Item(12,4).Value = 'FAPs'
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 2.489778666666667
$ws.Cells.Item(12,8).Value = 7.469336
$ws.Cells.Item(12,9).Value = 0.0008121061793956991
$ws.Cells.Item(12,10).Value = 0.0008121061793956991
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 1.27306
$ws.Cells.Item(12,14).Value = 3.81918
$ws.Cells.Item(12,15).Value = 0.0183646118807784
$ws.Cells.Item(12,16).Value = 0.0183646118807784
$ws.Cells.Item(12,17).Value = 3.169637629386666
$ws.Cells.Item(12,18).Value = 28.52673866448
$ws.Cells.Item(12,19).Value = [double]"1.491401479058381E-05"
$ws.Cells.Item(12,20).Value = [double]"1.491401479058381E-05"

$ws.Cells.Item(13,1).Value = 'Resolving-Mac'
$ws.Cells.Item(13,2).Value = 'Dcn'
$ws.Cells.Item(13,3).Value = 'Tlr2'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 2.489778666666667
$ws.Cells.Item(13,8).Value = 7.469336
$ws.Cells.Item(13,9).Value = 0.0008121061793956991
$ws.Cells.Item(13,10).Value = 0.0008121061793956991
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 66.45050433333334
$ws.Cells.Item(13,14).Value = 199.351513
$ws.Cells.Item(13,15).Value = 0.9585861792560053
$ws.Cells.Item(13,16).Value = 0.9585861792560054
$ws.Cells.Item(13,17).Value = 165.4470480783743
$ws.Cells.Item(13,18).Value = 1489.023432705368
$ws.Cells.Item(13,19).Value = 0.0007784737596571152
$ws.Cells.Item(13,20).Value = 0.0007784737596571153

$ws.Rows("14:17").Delete()
